$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (year 2025) with refreshed metrics
$ws.Range("C6").Value = 396
$ws.Range("D6").Value = 306
$ws.Range("F6").Value = 66.95842450765865
$ws.Range("G6").Value = 22.72727272727273
$ws.Range("H6").Value = 77.27272727272727
